# Fix mistake in deltaHrbolatosti func: append an explanatory sentence
# right after the word "pozice." at the end of the paragraph that reads
# "Tento výpočet je porovnán s dosavadní nejlepší pozicí a v případě
# bodového zlepšení je uložena pozice."
#
# The new sentence must be its own run (<w:r>), sharing the same
# run-formatting as the "pozice." run (rStyle "Odkazjemn", complex-script
# font Courier New, smallCaps off, automatic color, 12pt/12pt-cs,
# Czech language).

$d = $word.ActiveDocument

# Locate the unique sentence ending in "pozice." so we don't confuse it
# with the other occurrence of the word "pozice" elsewhere in the doc.
# NOTE: Find.Execute mutates the Range it was called on in place to span
# the found text, so we must keep using this same $searchRng afterwards
# (a fresh $d.Content call would give an unrelated, whole-document range).
$searchRng = $d.Content
$found = $searchRng.Find.Execute(
    "je uložena pozice.",  # FindText
    $false,                # MatchCase
    $false,                # MatchWholeWord
    $false,                # MatchWildcards
    $false,                # MatchSoundsLike
    $false,                # MatchAllWordForms
    $true,                 # Forward
    1,                     # Wrap (wdFindContinue)
    $false,                # Format
    "",                    # ReplaceWith
    0                      # Replace (wdReplaceNone)
)

if (-not $found) {
    throw "Could not find target sentence ending in 'pozice.'"
}

# Narrow the found range down to just the trailing "pozice." (7 chars),
# which lives in its own run with the formatting we need to copy.
$periodEnd = $searchRng.End
$periodStart = $periodEnd - 7
$periodRng = $d.Range($periodStart, $periodEnd)

if ($periodRng.Text -ne "pozice.") {
    throw "Unexpected text at target range: [$($periodRng.Text)]"
}

$newSentence = " Tyto konstanty byly získány experimentální cestou, avšak je jisté, že existuje mnohem lepší varianta."

# Replace the "pozice." range with itself plus the new sentence as a
# second run, both carrying the exact run properties of the original
# "pozice." run. Using InsertXML on this non-empty range lets Word merge
# the replacement in place (no stray paragraph break) while fully
# honoring the <w:rPr> we specify - unlike Range.InsertAfter/Font.*
# property setters, which this host does not reliably persist.
$xmlFrag = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:rPr><w:rStyle w:val="Odkazjemn"/><w:rFonts w:cs="Courier New"/><w:smallCaps w:val="0"/><w:color w:val="auto"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="cs-CZ"/></w:rPr><w:t>pozice.</w:t></w:r>
<w:r><w:rPr><w:rStyle w:val="Odkazjemn"/><w:rFonts w:cs="Courier New"/><w:smallCaps w:val="0"/><w:color w:val="auto"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="cs-CZ"/></w:rPr><w:t xml:space="preserve">$newSentence</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$periodRng.InsertXML($xmlFrag)

Write-Host "Inserted sentence after 'pozice.' successfully."
